# "Fix XML rappel with caution"
#
# Row 2 is corrected to reference the "Direction régionale" / Casablanca
# contract with the new gross/tax/net amounts; row 3 becomes a blank spacer
# row that only carries the repeated totals; the remaining detail rows
# (4-8, including the old bottom totals row) are removed entirely, which
# also shrinks the sheet's used range down to A1:K3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Drop the now-obsolete detail + totals rows (4-8) first ---
$ws.Rows("4:8").Delete()

# --- Row 2: update entity/contract + amounts ---
$ws.Range("F2").Value = "Direction régionale"
$ws.Range("G2").Value = "901/CASABLANCA/AV1"
$ws.Range("I2").Value = 8500.01
$ws.Range("J2").Value = 850.01
$ws.Range("K2").Value = 7650

# --- Row 3: blank out the identity columns, keep the repeated totals ---
$ws.Range("A3:H3").Value = " "
$ws.Range("I3").Value = 8500.01
$ws.Range("J3").Value = 850.01
$ws.Range("K3").Value = 7650
